$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (Förändrad) from 45207 to 45208 for rows 2 through 27
$ws.Range("C2:C27").Value = 45208
